$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that was sitting right after
#    the lone "I" run near the top of the document. Word keeps this
#    bookmark pinned to the very last edit location, so once we make a
#    fresh edit further down, it will naturally be re-created there.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Add the new "VISUALIZER - " lead-in right before the paragraph
#    that describes the viewer component ("Viewer will take the
#    generated files ..."), mirroring the existing "DATA PARSER - "
#    lead-in used earlier in the same section. Use InsertXML so the
#    new run carries the exact same run formatting (Arial Unicode MS
#    east-asian/complex-script font, size 24/24) as its neighboring
#    run, instead of plain-text insertion (which would leave the run
#    with no explicit rPr at all).
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.ClearFormatting()
$found = $target.Find.Execute("Viewer will take the generated files", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertStart = $target.Start
    $insertRange = $d.Range($insertStart, $insertStart)

    $xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">VISUALIZER - </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $insertRange.InsertXML($xmlFragment)

    # ------------------------------------------------------------------
    # 3) Drop Word's last-edit-location bookmark right after the text we
    #    just typed (between "VISUALIZER - " and "Viewer will take ...").
    # ------------------------------------------------------------------
    $goBackPos = $insertStart + [string]"VISUALIZER - ".Length
    $goBackRange = $d.Range($goBackPos, $goBackPos)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
